$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextCell {
    param($cell, $value)
    $origStyle = $cell.Style
    $cell.NumberFormat = "@"
    $cell.Value = $value
    $cell.Style = $origStyle
}

Set-TextCell $ws.Range('D2') '65.617.21'
Set-TextCell $ws.Range('E2') '  -0.86%  '
Set-TextCell $ws.Range('D3') '3.442.99'
Set-TextCell $ws.Range('E3') '  -3.08%  '
Set-TextCell $ws.Range('E4') '  -0.06%  '
Set-TextCell $ws.Range('D5') '592.60'
Set-TextCell $ws.Range('E5') '  -1.89%  '
Set-TextCell $ws.Range('D6') '137.47'
Set-TextCell $ws.Range('E6') '  -6.20%  '
Set-TextCell $ws.Range('D7') '3.441.49'
Set-TextCell $ws.Range('E7') '  -3.08%  '
Set-TextCell $ws.Range('E8') '  -0.02%  '
Set-TextCell $ws.Range('D9') '0.501'
Set-TextCell $ws.Range('E9') '  +1.06%  '
Set-TextCell $ws.Range('D10') '7.33'
Set-TextCell $ws.Range('E10') '  -5.96%  '
Set-TextCell $ws.Range('E11') '  -8.20%  '
Set-TextCell $ws.Range('E12') '  -7.01%  '
Set-TextCell $ws.Range('D13') '4.021.55'
Set-TextCell $ws.Range('E13') '  -3.11%  '
Set-TextCell $ws.Range('E14') '  -9.42%  '
Set-TextCell $ws.Range('D15') '26.48'
Set-TextCell $ws.Range('E15') '  -8.88%  '
Set-TextCell $ws.Range('D16') '3.434.88'
Set-TextCell $ws.Range('E16') '  -3.40%  '
Set-TextCell $ws.Range('D17') '65.575.55'
Set-TextCell $ws.Range('E17') '  -0.94%  '
Set-TextCell $ws.Range('E18') '  -1.99%  '
Set-TextCell $ws.Range('D19') '9.84'
Set-TextCell $ws.Range('E19') '  -10.88%  '
Set-TextCell $ws.Range('D20') '5.89'
Set-TextCell $ws.Range('E20') '  -5.38%  '
Set-TextCell $ws.Range('D21') '13.75'
Set-TextCell $ws.Range('E21') '  -6.51%  '
Set-TextCell $ws.Range('D22') '393.99'
Set-TextCell $ws.Range('E22') '  -5.79%  '
Set-TextCell $ws.Range('E23') '  -8.04%  '
Set-TextCell $ws.Range('D24') '73.39'
Set-TextCell $ws.Range('E24') '  -5.85%  '
Set-TextCell $ws.Range('E25') '  +0.02%  '
Set-TextCell $ws.Range('D26') '3.583.61'
Set-TextCell $ws.Range('E26') '  -3.00%  '
Set-TextCell $ws.Range('E27') '  -8.26%  '
Set-TextCell $ws.Range('E28') '  +0.12%  '
$ws.Range('B29').Value = 'InternetComputer(DFINITY)'
$ws.Range('C29').Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
Set-TextCell $ws.Range('D29') '8.26'
Set-TextCell $ws.Range('E29') '  -9.43%  '
$ws.Range('B30').Value = 'RenderToken'
$ws.Range('C30').Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
Set-TextCell $ws.Range('D30') '7.19'
Set-TextCell $ws.Range('E30') '  -8.71%  '
$ws.Range('B31').Value = 'PancakeSwap'
$ws.Range('C31').Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
Set-TextCell $ws.Range('D31') '2.25'
Set-TextCell $ws.Range('E31') '  -9.16%  '
Set-TextCell $ws.Range('D32') '3.448.12'
Set-TextCell $ws.Range('E32') '  -2.85%  '
Set-TextCell $ws.Range('E33') '  +0.02%  '
Set-TextCell $ws.Range('D34') '0.146'
Set-TextCell $ws.Range('E34') '  -6.09%  '
Set-TextCell $ws.Range('E35') '  -6.16%  '
Set-TextCell $ws.Range('D36') '172.76'
Set-TextCell $ws.Range('E36') '  -0.59%  '
Set-TextCell $ws.Range('D37') '6.95'
Set-TextCell $ws.Range('E37') '  -8.65%  '
Set-TextCell $ws.Range('D38') '1.18'
Set-TextCell $ws.Range('E38') '  -9.94%  '
Set-TextCell $ws.Range('E39') '  -7.60%  '
Set-TextCell $ws.Range('E40') '  -8.85%  '
Set-TextCell $ws.Range('D41') '0.0769'
Set-TextCell $ws.Range('E41') '  -6.55%  '
Set-TextCell $ws.Range('E42') '  -4.25%  '
Set-TextCell $ws.Range('D43') '43.78'
Set-TextCell $ws.Range('E43') '  -4.00%  '
Set-TextCell $ws.Range('E44') '  +0.01%  '
Set-TextCell $ws.Range('D45') '4.43'
Set-TextCell $ws.Range('E45') '  -13.07%  '
Set-TextCell $ws.Range('E46') '  -10.40%  '
Set-TextCell $ws.Range('D47') '23.07'
Set-TextCell $ws.Range('E47') '  +1.43%  '
Set-TextCell $ws.Range('D48') '1.11'
Set-TextCell $ws.Range('E48') '  +0.32%  '
Set-TextCell $ws.Range('D49') '6.59'
Set-TextCell $ws.Range('E49') '  -6.98%  '
Set-TextCell $ws.Range('D50') '2.10'
Set-TextCell $ws.Range('E50') '  -13.81%  '
Set-TextCell $ws.Range('D51') '2.213.71'
Set-TextCell $ws.Range('E51') '  -7.13%  '
